$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 813 (everything from 813 downward shifts down by 2).
$ws.Rows.Item(813).Insert()
$ws.Rows.Item(813).Insert()

# New row 813: Acelga "Primera" quality record for 2022-07-27 (serial 44769).
$ws.Range("A813").Value = 6
$ws.Range("B813").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C813").Value = "Metropolitana"
$ws.Range("D813").Value = 44769
$ws.Range("E813").Value = 13
$ws.Range("F813").Value = 100112009
$ws.Range("G813").Value = "Acelga"
$ws.Range("H813").Value = "Sin especificar"
$ws.Range("I813").Value = "Primera"
$ws.Range("J813").Value = 150
$ws.Range("K813").Value = 16000
$ws.Range("L813").Value = 16000
$ws.Range("M813").Value = 16000
$ws.Range("N813").Value = "$/docena de atados"
$ws.Range("O813").Value = "Región Metropolitana"
$ws.Range("P813").Value = 5333
$ws.Range("Q813").Value = 3
$ws.Range("R813").Value = "Hortaliza"

# New row 814: Acelga "Segunda" quality record for 2022-07-27 (serial 44769).
$ws.Range("A814").Value = 6
$ws.Range("B814").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C814").Value = "Metropolitana"
$ws.Range("D814").Value = 44769
$ws.Range("E814").Value = 13
$ws.Range("F814").Value = 100112009
$ws.Range("G814").Value = "Acelga"
$ws.Range("H814").Value = "Sin especificar"
$ws.Range("I814").Value = "Segunda"
$ws.Range("J814").Value = 90
$ws.Range("K814").Value = 12000
$ws.Range("L814").Value = 12000
$ws.Range("M814").Value = 12000
$ws.Range("N814").Value = "$/docena de atados"
$ws.Range("O814").Value = "Región Metropolitana"
$ws.Range("P814").Value = 4000
$ws.Range("Q814").Value = 3
$ws.Range("R814").Value = "Hortaliza"
